# Updates the "Jogos da Semana" odds sheet:
#  1) Removes the Atl. Nacional vs Ind. Medellin match (old row 9) entirely,
#     shifting the following matches (old rows 10-12) up by one row.
#  2) Refreshes a handful of odds values on the remaining, unmoved rows
#     (rows 2, 3, 7 and 8) to their latest quoted prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Delete the whole row for the Atl. Nacional x Ind. Medellin game ---
$ws.Rows(9).Delete()

# --- 2) Update odds on row 2 (Lanus x Platense) ---
$ws.Range("M2").Value = 1.14
$ws.Range("N2").Value = 5.5
$ws.Range("T2").Value = 2.08

# --- 3) Update odds on row 3 (Tigre x Defensa y Justicia) ---
$ws.Range("H3").Value = 3.2
$ws.Range("L3").Value = 4
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 2.75
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.62
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 2.5
$ws.Range("U3").Value = 1.91
$ws.Range("V3").Value = 1.91
$ws.Range("W3").Value = 7
$ws.Range("X3").Value = 10
$ws.Range("AA3").Value = 21
$ws.Range("AB3").Value = 34
$ws.Range("AC3").Value = 8
$ws.Range("AG3").Value = 9
$ws.Range("AK3").Value = 29
$ws.Range("AL3").Value = 41
$ws.Range("AM3").Value = 301
$ws.Range("AP3").Value = 26
$ws.Range("AS3").Value = 201
$ws.Range("AT3").Value = 2.5
$ws.Range("AU3").Value = 8.5
$ws.Range("AZ3").Value = 67
$ws.Range("BB3").Value = 251

# --- 4) Update odds on row 7 (Coritiba x Santos) ---
$ws.Range("M7").Value = 1.11
$ws.Range("N7").Value = 6.5

# --- 5) Update odds on row 8 (Paysandu PA x Brusque) ---
$ws.Range("O8").Value = 1.33
$ws.Range("P8").Value = 3.25
